# Add "Day 2" wrap-up content: a bold/colored heading-like line followed by
# two plain paragraphs, appended right after the last paragraph of the
# document (before the final section break).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the anchor paragraph - the last paragraph in the main story,
# which currently ends the document right before the sectPr.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs.Last

# ------------------------------------------------------------------
# Paragraph 1: "Planning Backend for this feature."
#   Large, bold, themed heading-like text - matches the look of the
#   other "Day N ..." headings (Heading2's resolved font/color) but
#   applied as direct character formatting rather than a paragraph
#   style reference.
# ------------------------------------------------------------------
$anchor.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Planning Backend for this feature."

$headingRange = $headingPara.Range
$headingRange.Font.NameAscii = "Calibri Light"
$headingRange.Font.NameFarEast = "Calibri Light"
$headingRange.Font.Bold = 1
$headingRange.Font.Size = 20
$headingRange.Font.TextColor.ObjectThemeColor = 4
$headingRange.Font.TextColor.RGB = 0x96542F

# ------------------------------------------------------------------
# Paragraph 2: "Gemman Kavach Vision will be a seprate backend, ..."
# ------------------------------------------------------------------
$headingPara.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Gemman"
$p2End = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$p2End.InsertAfter(" Kavach Vision will be a seprate backend, running so some port, the ui will seprate so we do design a full backend gemma kavach vision.")

# ------------------------------------------------------------------
# Paragraph 3: "Also we need a object storage location for now using ..."
# ------------------------------------------------------------------
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Also we need a object storage location for now using "
$p3End = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$p3End.InsertAfter(" gemma3n-raw")
$p3End2 = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$p3End2.InsertAfter(" a goole coud storage engine")

Write-Output "Inserted Day 2 wrap-up paragraphs."
